$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.099.57'
$ws.Range('E2').Value = '  -0.51%  '

$ws.Range('D3').Value = '1.836.84'
$ws.Range('E3').Value = '  -0.42%  '

$ws.Range('E4').Value = '  +0.41%  '

$ws.Range('D5').Value = "'" + '242.66'
$ws.Range('E5').Value = '  -0.05%  '

$ws.Range('D6').Value = "'" + '0.6240'
$ws.Range('E6').Value = '  -5.83%  '

$ws.Range('D7').Value = "'" + '1.004'
$ws.Range('E7').Value = '  +0.41%  '

$ws.Range('D8').Value = "'" + '0.07572'
$ws.Range('E8').Value = '  +1.76%  '

$ws.Range('D9').Value = "'" + '0.2919'
$ws.Range('E9').Value = '  -1.28%  '

$ws.Range('D10').Value = "'" + '22.64'
$ws.Range('E10').Value = '  -2.90%  '

$ws.Range('D11').Value = "'" + '0.07765'
$ws.Range('E11').Value = '  +0.03%  '

$ws.Range('D12').Value = '1.836.69'
$ws.Range('E12').Value = '  -0.40%  '

$ws.Range('D13').Value = "'" + '4.951'
$ws.Range('E13').Value = '  -1.48%  '

$ws.Range('D14').Value = "'" + '0.6652'
$ws.Range('E14').Value = '  -1.15%  '

$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'" + '82.66'
$ws.Range('E15').Value = '  -0.91%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = "'" + '0.000009923'
$ws.Range('E16').Value = '  +13.85%  '

$ws.Range('D17').Value = "'" + '6.028'

$ws.Range('D18').Value = '29.152.68'
$ws.Range('E18').Value = '  -0.32%  '

$ws.Range('D19').Value = "'" + '225.32'
$ws.Range('E19').Value = '  -0.69%  '

$ws.Range('D20').Value = "'" + '12.34'
$ws.Range('E20').Value = '  -1.66%  '

$ws.Range('D21').Value = "'" + '1.004'
$ws.Range('E21').Value = '  +0.33%  '

$ws.Range('D22').Value = "'" + '7.196'
$ws.Range('E22').Value = '  +0.47%  '

$ws.Range('E23').Value = '  +0.45%  '

$ws.Range('D24').Value = "'" + '158.96'
$ws.Range('E24').Value = '  +0.04%  '

$ws.Range('D25').Value = "'" + '8.459'
$ws.Range('E25').Value = '  -2.01%  '

$ws.Range('D26').Value = "'" + '0.1365'
$ws.Range('E26').Value = '  -2.93%  '

$ws.Range('D27').Value = "'" + '17.90'
$ws.Range('E27').Value = '  -0.77%  '

$ws.Range('D28').Value = "'" + '1.495'
$ws.Range('E28').Value = '  -0.93%  '

$ws.Range('D29').Value = "'" + '4.072'
$ws.Range('E29').Value = '  -1.75%  '

$ws.Range('D30').Value = "'" + '4.037'
$ws.Range('E30').Value = '  -0.66%  '

$ws.Range('E31').Value = '  +0.91%  '

$ws.Range('D32').Value = "'" + '0.05203'
$ws.Range('E32').Value = '  -2.40%  '

$ws.Range('D33').Value = "'" + '1.850'
$ws.Range('E33').Value = '  -1.00%  '

$ws.Range('D34').Value = "'" + '0.7374'
$ws.Range('E34').Value = '  -1.49%  '

$ws.Range('D35').Value = "'" + '1.142'
$ws.Range('E35').Value = '  -1.42%  '

$ws.Range('D36').Value = "'" + '2.708'
$ws.Range('E36').Value = '  +2.01%  '

$ws.Range('D37').Value = '1.252.47'
$ws.Range('E37').Value = '  -4.68%  '

$ws.Range('D38').Value = "'" + '2.766'
$ws.Range('E38').Value = '  +0.30%  '

$ws.Range('D39').Value = "'" + '0.01782'
$ws.Range('E39').Value = '  -1.12%  '

$ws.Range('D40').Value = "'" + '6.318'
$ws.Range('E40').Value = '  -1.15%  '

$ws.Range('D41').Value = "'" + '0.8971'
$ws.Range('E41').Value = '  -1.08%  '

$ws.Range('D42').Value = "'" + '1.005'
$ws.Range('E42').Value = '  +0.49%  '

$ws.Range('D43').Value = "'" + '101.46'
$ws.Range('E43').Value = '  -2.12%  '

$ws.Range('D44').Value = '1.979.49'
$ws.Range('E44').Value = '  -0.55%  '

$ws.Range('D45').Value = "'" + '0.00000000125'
$ws.Range('E45').Value = '  +2.14%  '

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = "'" + '0.5127'
$ws.Range('E46').Value = '  -0.32%  '

$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = "'" + '64.03'
$ws.Range('E47').Value = '  -1.95%  '

$ws.Range('D48').Value = "'" + '0.4005'
$ws.Range('E48').Value = '  -0.45%  '

$ws.Range('D49').Value = "'" + '8.823'
$ws.Range('E49').Value = '  +0.65%  '

$ws.Range('D50').Value = "'" + '0.05761'
$ws.Range('E50').Value = '  -1.88%  '

$ws.Range('E51').Value = '  -6.31%  '

